$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.173520803451538
$ws.Range("B1").Value = 3.058849096298218
$ws.Range("C1").Value = 2.632047653198242
$ws.Range("D1").Value = 2.160658597946167
$ws.Range("E1").Value = 1.454070329666138
